# Laborator 28.03.2023 - marking "sapt 5" (week 5, column G) attendance
# as present (TRUE) for a set of students.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows in which column G (week 5 attendance) should be set to TRUE
$rows = @(3, 5, 7, 14, 15, 18, 21, 23, 24, 26, 30, 36, 37, 38)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 7).Value = $true
}

# Update the view: scroll so row 22 is the top-left visible row,
# and move the active selection to L36.
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("L36").Select()
